$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 date moved ahead by one day (45308 -> 45309, i.e. 2024-01-17 -> 2024-01-18)
$ws.Range("A1").Value = 45309

# Price list "step 1 and 2" increase (~53%) for rows 33-38 (column D)
$ws.Range("D33").Value = 3823.437
$ws.Range("D34").Value = 3823.437
$ws.Range("D35").Value = 5001.688
$ws.Range("D36").Value = 5001.725
$ws.Range("D37").Value = 5001.688
$ws.Range("D38").Value = 6304.761
